$wb = $excel.ActiveWorkbook

# --- Sheet "2025-2" (sheet1): insert a new registry row for "EP TASA 35" ---
$ws1 = $wb.Worksheets.Item("2025-2")

# A new boat ("EP TASA 35" / GP/87) was registered and inserted right after
# the "EP TASA 34" row (row 5), pushing the remolcador/EP/remolcador/EP rows
# that followed down by one.
$ws1.Rows.Item(6).Insert()

$ws1.Range("A6").Value = "2025-2"
$ws1.Range("B6").Value = "EP TASA 35"
$ws1.Range("C6").Value = "Embarcación Pesquera"
$ws1.Range("D6").Value = "GP/87"
$ws1.Range("E6").Value = "GP/87-225"

# --- Sheet "desguace" (sheet2): tidy up the duplicate/stray style on C5:C6 ---
$ws2 = $wb.Worksheets.Item("desguace")

# C5 and C6 were carrying a redundant cell format (a near-duplicate of the
# "Nave" column style used everywhere else, just with an unnecessary
# applyFill flag). Repoint them to the same format already used by the
# other "Nave" cells (e.g. C2) instead of that stray one.
$ws2.Range("C2").Copy()
$ws2.Range("C5:C6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection / active-sheet bookkeeping left by the editing session ---
$ws2.Activate()
$ws2.Range("B9:C9").Select()

$ws1.Activate()
$ws1.Range("G9").Select()
